# Update "want to go" counts (column F) for several rows across sheets.
# Sheet 1 = 展览 (Exhibition), Sheet 3 = 本地生活 (Local Life), Sheet 4 = 全部类型 (All Types)

$wb = $excel.ActiveWorkbook

# --- Sheet 1: 展览 ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("F2").Value  = 550
$ws1.Range("F5").Value  = 1365
$ws1.Range("F6").Value  = 674
$ws1.Range("F7").Value  = 359
$ws1.Range("F9").Value  = 159
$ws1.Range("F10").Value = 429
$ws1.Range("F11").Value = 6337
$ws1.Range("F14").Value = 1896
$ws1.Range("F15").Value = 4740
$ws1.Range("F16").Value = 468
$ws1.Range("F19").Value = 5521
$ws1.Range("F20").Value = 7309
$ws1.Range("F30").Value = 1502
$ws1.Range("F31").Value = 569
$ws1.Range("F32").Value = 697
$ws1.Range("F42").Value = 1443
$ws1.Range("F49").Value = 3963

# --- Sheet 3: 本地生活 ---
$ws3 = $wb.Worksheets.Item(3)
$ws3.Range("F2").Value = 4466

# --- Sheet 4: 全部类型 ---
$ws4 = $wb.Worksheets.Item(4)
$ws4.Range("F2").Value  = 4466
$ws4.Range("F3").Value  = 550
$ws4.Range("F9").Value  = 1365
$ws4.Range("F11").Value = 674
$ws4.Range("F12").Value = 359
$ws4.Range("F13").Value = 159
$ws4.Range("F14").Value = 429
$ws4.Range("F17").Value = 4740
$ws4.Range("F18").Value = 5521
$ws4.Range("F19").Value = 5521
$ws4.Range("F29").Value = 1502
$ws4.Range("F30").Value = 569
$ws4.Range("F31").Value = 697
$ws4.Range("F50").Value = 3963
